# Adds the "Vendredi 17 mars" journal entry, replacing the two trailing
# empty paragraphs at the end of the document with the new heading /
# sub-heading / body paragraphs described in the commit message
# ("Correction de plusieurs bugs et optimisation").

$d = $word.ActiveDocument
$cr = [char]13

# Locate the two trailing empty paragraphs (each containing nothing but
# a paragraph mark) that sit right before the final section break.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$prevPara = $d.Paragraphs.Item($count - 1)

if ($lastPara.Range.Text -eq $cr -and $prevPara.Range.Text -eq $cr) {
    # Drop the first of the two empty paragraphs entirely so only one
    # (empty) paragraph remains at the tail of the document.
    $null = $prevPara.Range.Delete()
}

# Re-fetch the now-last paragraph; InsertXML replaces the paragraph
# targeted by the range, so this single call drops the remaining empty
# paragraph and substitutes the new heading/sub-heading/body paragraphs.
$count = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count)

$newContentXml = @"
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Titre1"/>
      </w:pPr>
      <w:r>
        <w:t>Vendredi 17 mars :</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="Titre2"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Affichage des données : </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:jc w:val="both"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Je reprends donc le projet en fin de matinée vers les 11h </w:t>
      </w:r>
      <w:r>
        <w:t>dû</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> à encore une fois des problèmes techniques. Je </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">reprends donc l’affichage des données entamé </w:t>
      </w:r>
      <w:r>
        <w:t>la vieille. Tout en réglant plein de petites choses à droites. Notamment j’ai ajouté deux fonctions l’une qui consiste à formater la date afin de la rendre lisible à l’utilisateur et une qui consiste à vérifier si un contenu est marqué comme explicit et si tel est le cas à rajouter à la fin de son titre un petit « </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:cs="MS Gothic" w:hint="eastAsia"/>
        </w:rPr>
        <w:t>ⓔ</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t> »</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>. De plus j’ai ajouté dans le tri de mes données une catégorie music-</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>video</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> car jusqu’à présent je ne </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>l’avais</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> pas vu mais une requête de type musique pouvait renvoyer un clip et ce type étant catégorisé comme default avait pour conséquence d’afficher une </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>card</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> vide. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>A ce point j’ai push mon avancement sous le nom « Correction de plusieurs bugs et optimisation »</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> vers les 12h</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">Une fois fais-je commence le vrai nom de cette catégorie le paramétrage de l’affichage des données. Une grande partie du traitement ce fait en </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>golang</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve"> directement avec ma nouvelle structure, mais il est toujours nécessaire de faire des ajustements avec des if sur le </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t>template</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="MS Gothic" w:cstheme="minorHAnsi"/>
        </w:rPr>
        <w:t xml:space="preserve">, comme pour par exemple choisir entre afficher une description ou un album en fonction d’un livre ou d’une musique. </w:t>
      </w:r>
    </w:p>
"@

$null = $target.Range.InsertXML($newContentXml)
